$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column (C) for all data rows
# from serial date 45204 to 45207, leaving everything else untouched.
$ws.Range("C2:C62").Value = 45207
